$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 39, pushing the existing rows 39-157 down to 40-158.
$ws.Rows("39:39").Insert()

# Populate the newly inserted row 39 with the new price record.
$ws.Range("A39").Value = 4
$ws.Range("B39").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C39").Value = "Los Lagos"
$ws.Range("D39").Value = 45076
$ws.Range("E39").Value = 10
$ws.Range("F39").Value = "Fruta"
$ws.Range("G39").Value = 100104
$ws.Range("H39").Value = "Frutos de pepita"
$ws.Range("I39").Value = 100104003
$ws.Range("J39").Value = "Membrillo"
$ws.Range("K39").Value = "Champion"
$ws.Range("L39").Value = "Primera"
$ws.Range("M39").Value = 200
$ws.Range("N39").Value = 13000
$ws.Range("O39").Value = 14000
$ws.Range("P39").Value = 13500
$ws.Range("Q39").Value = "`$/caja 18 kilos empedrada"
$ws.Range("R39").Value = "Región de O'Higgins"
$ws.Range("S39").Value = 750
$ws.Range("T39").Value = 18
